$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 15680
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
